# Auto-generated Excel COM-interop edit script
# Applies header renames, Spanish title-case fixes to municipality names,
# floating point precision corrections, and removes trailing footer/notes rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case field names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Body cells: title-case Spanish prepositions/articles (de/del/la/las/el/los/y) ---
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("B21").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Bejucal De Ocampo"
$ws.Range("B31").Value = "Comitán De Domínguez"
$ws.Range("B47").Value = "Mazapa De Madero"
$ws.Range("B55").Value = "Salto De Agua"
$ws.Range("B56").Value = "San Cristóbal De Las Casas"
$ws.Range("B82").Value = "Hidalgo Del Parral"
$ws.Range("A98").Value = "Ciudad De México"
$ws.Range("B102").Value = "Cuajimalpa De Morelos"
$ws.Range("B124").Value = "San Juan Del Río"
$ws.Range("A127").Value = "Estado De México"
$ws.Range("B127").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B130").Value = "Almoloya De Alquisiras"
$ws.Range("B131").Value = "Almoloya De Juárez"
$ws.Range("B132").Value = "Almoloya Del Río"
$ws.Range("B137").Value = "Atizapán De Zaragoza"
$ws.Range("B142").Value = "Chapa De Mota"
$ws.Range("B146").Value = "Coacalco De Berriozábal"
$ws.Range("B152").Value = "Ecatepec De Morelos"
$ws.Range("B158").Value = "Ixtapan De La Sal"
$ws.Range("B159").Value = "Ixtapan Del Oro"
$ws.Range("B172").Value = "Naucalpan De Juárez"
$ws.Range("B182").Value = "San Felipe Del Progreso"
$ws.Range("B183").Value = "San Martín De Las Pirámides"
$ws.Range("B195").Value = "Tenango Del Valle"
$ws.Range("B206").Value = "Tlalnepantla De Baz"
$ws.Range("B212").Value = "Valle De Bravo"
$ws.Range("B213").Value = "Valle De Chalco Solidaridad"
$ws.Range("B214").Value = "Villa De Allende"
$ws.Range("B215").Value = "Villa Del Carbón"
$ws.Range("B225").Value = "San Miguel De Allende"
$ws.Range("B226").Value = "Apaseo El Alto"
$ws.Range("B227").Value = "Apaseo El Grande"
$ws.Range("B235").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B248").Value = "San Luis De La Paz"
$ws.Range("B249").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B251").Value = "Silao De La Victoria"
$ws.Range("B254").Value = "Valle De Santiago"
$ws.Range("B258").Value = "Acapulco De Juárez"
$ws.Range("B261").Value = "Ajuchitlán Del Progreso"
$ws.Range("B262").Value = "Alcozauca De Guerrero"
$ws.Range("B266").Value = "Atenango Del Río"
$ws.Range("B267").Value = "Atlamajalcingo Del Monte"
$ws.Range("B269").Value = "Atoyac De Álvarez"
$ws.Range("B270").Value = "Ayutla De Los Libres"
$ws.Range("B273").Value = "Buenavista De Cuéllar"
$ws.Range("B274").Value = "Chilapa De Álvarez"
$ws.Range("B275").Value = "Chilpancingo De Los Bravo"
$ws.Range("B276").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B281").Value = "Coyuca De Benítez"
$ws.Range("B282").Value = "Coyuca De Catalán"
$ws.Range("B286").Value = "Cuetzala Del Progreso"
$ws.Range("B287").Value = "Cutzamala De Pinzón"
$ws.Range("B293").Value = "Huitzuco De Los Figueroa"
$ws.Range("B294").Value = "Iguala De La Independencia"
$ws.Range("B296").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B297").Value = "Zihuatanejo De Azueta"
$ws.Range("B299").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B302").Value = "Mártir De Cuilapan"
$ws.Range("B315").Value = "Taxco De Alarcón"
$ws.Range("B317").Value = "Técpan De Galeana"
$ws.Range("B319").Value = "Tepecoacuilco De Trujano"
$ws.Range("B321").Value = "Tixtla De Guerrero"
$ws.Range("B325").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B326").Value = "Tlapa De Comonfort"
$ws.Range("B337").Value = "Agua Blanca De Iturbide"
$ws.Range("B342").Value = "Atotonilco El Grande"
$ws.Range("B348").Value = "Cuautepec De Hinojosa"
$ws.Range("B351").Value = "Huasca De Ocampo"
$ws.Range("B355").Value = "Huejutla De Reyes"
$ws.Range("B358").Value = "Jacala De Ledezma"
$ws.Range("B364").Value = "Mineral Del Chico"
$ws.Range("B365").Value = "Mineral Del Monte"
$ws.Range("B366").Value = "Mixquiahuala De Juárez"
$ws.Range("B368").Value = "Nopala De Villagrán"
$ws.Range("B369").Value = "Omitlán De Juárez"
$ws.Range("B370").Value = "Pachuca De Soto"
$ws.Range("B373").Value = "Progreso De Obregón"
$ws.Range("B378").Value = "Santiago De Anaya"
$ws.Range("B379").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B383").Value = "Tenango De Doria"
$ws.Range("B385").Value = "Tepehuacán De Guerrero"
$ws.Range("B386").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B389").Value = "Tezontepec De Aldama"
$ws.Range("B397").Value = "Tula De Allende"
$ws.Range("B398").Value = "Tulancingo De Bravo"
$ws.Range("B399").Value = "Villa De Tezontepec"
$ws.Range("B403").Value = "Zacualtipán De Ángeles"
$ws.Range("B409").Value = "Autlán De Navarro"
$ws.Range("B417").Value = "Huejuquilla El Alto"
$ws.Range("B422").Value = "La Manzanilla De La Paz"
$ws.Range("B423").Value = "Lagos De Moreno"
$ws.Range("B431").Value = "San Juan De Los Lagos"
$ws.Range("B434").Value = "Tamazula De Gordiano"
$ws.Range("B438").Value = "Tizapán El Alto"
$ws.Range("B439").Value = "Tlajomulco De Zúñiga"
$ws.Range("B443").Value = "Unión De San Antonio"
$ws.Range("B444").Value = "Unión De Tula"
$ws.Range("B445").Value = "Valle De Juárez"
$ws.Range("B447").Value = "Yahualica De González Gallo"
$ws.Range("B465").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B515").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B543").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B547").Value = "Puente De Ixtla"
$ws.Range("B553").Value = "Tetela Del Volcán"
$ws.Range("B555").Value = "Tlaltizapán De Zapata"
$ws.Range("B563").Value = "Zacualpan De Amilpas"
$ws.Range("B580").Value = "San Nicolás De Los Garza"
$ws.Range("B584").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B589").Value = "Ayoquezco De Aldama"
$ws.Range("B594").Value = "Chalcatongo De Hidalgo"
$ws.Range("B595").Value = "Chiquihuitlán De Benito Juárez"
$ws.Range("B596").Value = "Ciénega De Zimatlán"
$ws.Range("B599").Value = "Coicoyán De Las Flores"
$ws.Range("B600").Value = "Constancia Del Rosario"
$ws.Range("B603").Value = "Cuilápam De Guerrero"
$ws.Range("B604").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B605").Value = "El Barrio De La Soledad"
$ws.Range("B606").Value = "Fresnillo De Trujano"
$ws.Range("B607").Value = "Guadalupe De Ramírez"
$ws.Range("B608").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B609").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B610").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B611").Value = "Huautla De Jiménez"
$ws.Range("B613").Value = "Ixtlán De Juárez"
$ws.Range("B614").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B617").Value = "Mariscala De Juárez"
$ws.Range("B619").Value = "Mazatlán Villa De Flores"
$ws.Range("B621").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B622").Value = "Mixistlán De La Reforma"
$ws.Range("B625").Value = "Nejapa De Madero"
$ws.Range("B627").Value = "Oaxaca De Juárez"
$ws.Range("B628").Value = "Ocotlán De Morelos"
$ws.Range("B629").Value = "Pinotepa De Don Luis"
$ws.Range("B631").Value = "Putla Villa De Guerrero"
$ws.Range("B632").Value = "Reforma De Pineda"
$ws.Range("B634").Value = "Rojas De Cuauhtémoc"
$ws.Range("B649").Value = "San Antonio De La Cal"
$ws.Range("B677").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B685").Value = "San Juan De Los Cués"
$ws.Range("B686").Value = "San Juan Del Estado"
$ws.Range("B687").Value = "San Juan Del Río"
$ws.Range("B726").Value = "San Miguel Del Puerto"
$ws.Range("B743").Value = "San Pablo Villa De Mitla"
$ws.Range("B762").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B774").Value = "Santa Ana Del Valle"
$ws.Range("B783").Value = "Santa Cruz De Bravo"
$ws.Range("B788").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B791").Value = "Santa Inés De Zaragoza"
$ws.Range("B792").Value = "Santa Inés Del Monte"
$ws.Range("B793").Value = "Santa Lucía Del Camino"
$ws.Range("B805").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B852").Value = "Santo Domingo De Morelos"
$ws.Range("B865").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B866").Value = "Tanetze De Zaragoza"
$ws.Range("B868").Value = "Tataltepec De Valdés"
$ws.Range("B869").Value = "Teotitlán De Flores Magón"
$ws.Range("B870").Value = "Teotitlán Del Valle"
$ws.Range("B872").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B873").Value = "Tlacolula De Matamoros"
$ws.Range("B874").Value = "Tlalixtac De Cabrera"
$ws.Range("B875").Value = "Totontepec Villa De Morelos"
$ws.Range("B878").Value = "Villa De Chilapa De Díaz"
$ws.Range("B879").Value = "Villa De Etla"
$ws.Range("B880").Value = "Villa De Tututepec"
$ws.Range("B881").Value = "Villa De Zaachila"
$ws.Range("B884").Value = "Villa Sola De Vega"
$ws.Range("B885").Value = "Yutanduchi De Guerrero"
$ws.Range("B888").Value = "Zimatlán De Álvarez"
$ws.Range("B914").Value = "Ayotoxco De Guerrero"
$ws.Range("B920").Value = "Chalchicomula De Sesma"
$ws.Range("B931").Value = "Chila De La Sal"
$ws.Range("B942").Value = "Cuapiaxtla De Madero"
$ws.Range("B946").Value = "Cuayuca De Andrade"
$ws.Range("B947").Value = "Cuetzalan Del Progreso"
$ws.Range("B963").Value = "Huehuetlán El Chico"
$ws.Range("B964").Value = "Huehuetlán El Grande"
$ws.Range("B968").Value = "Huitzilan De Serdán"
$ws.Range("B970").Value = "Ixcamilpa De Guerrero"
$ws.Range("B973").Value = "Izúcar De Matamoros"
$ws.Range("B984").Value = "Los Reyes De Juárez"
$ws.Range("B985").Value = "Mazapiltepec De Juárez"
$ws.Range("B998").Value = "Palmar De Bravo"
$ws.Range("B1008").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B1025").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1029").Value = "San Salvador El Seco"
$ws.Range("B1030").Value = "San Salvador El Verde"
$ws.Range("B1039").Value = "Tecali De Herrera"
$ws.Range("B1047").Value = "Tepanco De López"
$ws.Range("B1048").Value = "Tepango De Rodríguez"
$ws.Range("B1049").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1055").Value = "Tepexi De Rodríguez"
$ws.Range("B1057").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B1058").Value = "Tetela De Ocampo"
$ws.Range("B1059").Value = "Teteles De Avila Castillo"
$ws.Range("B1064").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1076").Value = "Totoltepec De Guerrero"
$ws.Range("B1078").Value = "Tuzamapan De Galeana"
$ws.Range("B1082").Value = "Xayacatlán De Bravo"
$ws.Range("B1088").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1103").Value = "Amealco De Bonfil"
$ws.Range("B1104").Value = "Cadereyta De Montes"
$ws.Range("B1106").Value = "Jalpan De Serra"
$ws.Range("B1107").Value = "Landa De Matamoros"
$ws.Range("B1109").Value = "Pinal De Amoles"
$ws.Range("B1111").Value = "San Juan Del Río"
$ws.Range("B1120").Value = "Ciudad Del Maíz"
$ws.Range("B1126").Value = "Mexquitic De Carmona"
$ws.Range("B1135").Value = "Santa María Del Río"
$ws.Range("B1142").Value = "Tanquián De Escobedo"
$ws.Range("B1144").Value = "Villa De Arista"
$ws.Range("B1145").Value = "Villa De Guadalupe"
$ws.Range("B1146").Value = "Villa De Reyes"
$ws.Range("B1165").Value = "Nacozari De García"
$ws.Range("B1192").Value = "Soto La Marina"
$ws.Range("B1198").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1200").Value = "Amaxac De Guerrero"
$ws.Range("B1201").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B1206").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1212").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1215").Value = "Mazatecochco De José María Morelos"
$ws.Range("B1216").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1219").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1222").Value = "San Pablo Del Monte"
$ws.Range("B1223").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B1231").Value = "Tepetitla De Lardizábal"
$ws.Range("B1234").Value = "Tetla De La Solidaridad"
$ws.Range("B1255").Value = "Amatlán De Los Reyes"
$ws.Range("B1262").Value = "Boca Del Río"
$ws.Range("B1267").Value = "Castillo De Teayo"
$ws.Range("B1269").Value = "Cazones De Herrera"
$ws.Range("B1280").Value = "Cosamaloapan De Carpio"
$ws.Range("B1296").Value = "Hueyapan De Ocampo"
$ws.Range("B1297").Value = "Ignacio De La Llave"
$ws.Range("B1301").Value = "Ixhuatlán De Madero"
$ws.Range("B1302").Value = "Ixhuatlán Del Café"
$ws.Range("B1303").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1311").Value = "Juchique De Ferrer"
$ws.Range("B1314").Value = "Las Vigas De Ramírez"
$ws.Range("B1315").Value = "Lerdo De Tejada"
$ws.Range("B1319").Value = "Martínez De La Torre"
$ws.Range("B1321").Value = "Medellín De Bravo"
$ws.Range("B1325").Value = "Mixtla De Altamirano"
$ws.Range("B1334").Value = "Paso De Ovejas"
$ws.Range("B1335").Value = "Paso Del Macho"
$ws.Range("B1339").Value = "Poza Rica De Hidalgo"
$ws.Range("B1346").Value = "Sayula De Alemán"
$ws.Range("B1349").Value = "Soledad De Doblado"
$ws.Range("B1353").Value = "Tatahuicapan De Juárez"
$ws.Range("B1368").Value = "Tlacotepec De Mejía"
$ws.Range("B1377").Value = "Vega De Alatorre"
$ws.Range("B1386").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1387").Value = "Zozocolco De Hidalgo"
$ws.Range("B1401").Value = "Concepción Del Oro"

# --- Floating point precision corrections (pct_matriculas recomputed with 1 ULP diff) ---
$ws.Range("D102").Value = 0.0009044862518089724
$ws.Range("D381").Value = 0.0009044862518089724
$ws.Range("D451").Value = 0.009804630969609264
$ws.Range("D631").Value = 0.0009044862518089724
$ws.Range("D910").Value = 0.0009044862518089724
$ws.Range("D999").Value = 0.0009044862518089724
$ws.Range("D1092").Value = 0.0009044862518089724
$ws.Range("D1233").Value = 0.0009044862518089724
$ws.Range("D1343").Value = 0.0009044862518089724

# --- Remove trailing footer/notes rows (1415-1419); row 1414 was already empty ---
$ws.Range("A1415:A1419").EntireRow.Delete()
